$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -7.779999999999999
$ws.Range("C7").Value = -12.844
$ws.Range("A10").Value = -21.776
$ws.Range("A12").Value = -21.696
$ws.Range("C15").Value = -14.056
$ws.Range("A18").Value = -22.265
$ws.Range("D18").Value = -8.84
$ws.Range("D19").Value = -8.113
$ws.Range("C20").Value = -12.673
$ws.Range("D27").Value = -8.744999999999999
$ws.Range("C29").Value = -12.124
$ws.Range("C30").Value = -13.018
$ws.Range("C31").Value = -13.358
$ws.Range("A37").Value = -20.029
$ws.Range("C40").Value = -12.782
$ws.Range("D42").Value = -8.257999999999999
$ws.Range("D44").Value = -7.502
$ws.Range("D47").Value = -7.291999999999999
$ws.Range("A55").Value = -21.894
$ws.Range("D58").Value = -8.400999999999998
$ws.Range("A68").Value = -21.736
$ws.Range("C68").Value = -11.001
$ws.Range("D73").Value = -8.105
$ws.Range("C76").Value = -13.045
$ws.Range("A77").Value = -20.843
$ws.Range("A78").Value = -19.951
$ws.Range("C87").Value = -13.458
$ws.Range("C88").Value = -13.35
$ws.Range("D95").Value = -7.718999999999999
$ws.Range("C96").Value = -12.644
$ws.Range("C98").Value = -13.23
$ws.Range("C101").Value = -13.048
$ws.Range("D101").Value = -8.061999999999999
$ws.Range("C102").Value = -13.086
